$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'21.765.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.74%  "

# Row 3
$ws.Range("D3").Value = "'1.539.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.49%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("E5").Value = "  +0.05%  "

# Row 6
$ws.Range("D6").Value = "'290.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "

# Row 7
$ws.Range("D7").Value = "'0.3877"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.97%  "

# Row 8
$ws.Range("D8").Value = "'0.3191"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.10%  "

# Row 9
$ws.Range("D9").Value = "'43.15"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.53%  "

# Row 10
$ws.Range("D10").Value = "'0.07197"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.48%  "

# Row 11
$ws.Range("D11").Value = "'1.061"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.93%  "

# Row 12
$ws.Range("E12").Value = "  +0.07%  "

# Row 13
$ws.Range("D13").Value = "'5.637"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.31%  "

# Row 15
$ws.Range("D15").Value = "'6.628"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.05%  "

# Row 16
$ws.Range("D16").Value = "'1.543.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.86%  "

# Row 17
$ws.Range("D17").Value = "'0.00001109"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.11%  "

# Row 18
$ws.Range("D18").Value = "'0.06586"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.08%  "

# Row 19
$ws.Range("D19").Value = "'83.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.80%  "

# Row 20
$ws.Range("E20").Value = "  +0.08%  "

# Row 21
$ws.Range("D21").Value = "'6.145"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.83%  "

# Row 22
$ws.Range("D22").Value = "'15.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.76%  "

# Row 23
$ws.Range("D23").Value = "'10.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.61%  "

# Row 24
$ws.Range("E24").Value = "  +4.73%  "

# Row 25
$ws.Range("D25").Value = "'21.770.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.77%  "

# Row 26
$ws.Range("D26").Value = "'2.395"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.44%  "

# Row 27
$ws.Range("D27").Value = "'146.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.09%  "

# Row 28
$ws.Range("E28").Value = "  -3.95%  "

# Row 29
$ws.Range("D29").Value = "'4.850"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.34%  "

# Row 30
$ws.Range("D30").Value = "'1.718.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.59%  "

# Row 31
$ws.Range("D31").Value = "'117.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.22%  "

# Row 32
$ws.Range("D32").Value = "'0.9674"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -14.21%  "

# Row 33
$ws.Range("D33").Value = "'5.898"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.51%  "

# Row 34
$ws.Range("D34").Value = "'0.08215"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.36%  "

# Row 35
$ws.Range("D35").Value = "'8.928"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.07%  "

# Row 36
$ws.Range("D36").Value = "'0.06088"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.34%  "

# Row 37
$ws.Range("D37").Value = "'5.144"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.98%  "

# Row 38
$ws.Range("D38").Value = "'1.481"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -19.74%  "

# Row 39
$ws.Range("D39").Value = "'0.02208"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.09%  "

# Row 40
$ws.Range("E40").Value = "  -4.92%  "

# Row 41
$ws.Range("D41").Value = "'1.188"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.23%  "

# Row 42
$ws.Range("E42").Value = "  +0.03%  "

# Row 43
$ws.Range("D43").Value = "'10.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.63%  "

# Row 44
$ws.Range("D44").Value = "'0.5758"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.11%  "

# Row 45
$ws.Range("D45").Value = "'13.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.73%  "

# Row 46
$ws.Range("D46").Value = "'3.745"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.18%  "

# Row 47
$ws.Range("D47").Value = "'0.5519"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.91%  "

# Row 48
$ws.Range("D48").Value = "'118.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.24%  "

# Row 49
$ws.Range("D49").Value = "'1.866"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.84%  "

# Row 50
$ws.Range("D50").Value = "'1.144"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.37%  "

# Row 51
$ws.Range("D51").Value = "'0.06733"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.60%  "
